# Daily Report update: 2026-02-10 (data as of 2026-02-09, serial 46062)
# Adds a new day's rows to Daily_Data, and refreshes the Today_Summary
# and Monthly_Stats rollups that depend on the latest totals.

$wb = $excel.ActiveWorkbook
$wsDaily   = $wb.Worksheets.Item("Daily_Data")
$wsSummary = $wb.Worksheets.Item("Today_Summary")
$wsMonthly = $wb.Worksheets.Item("Monthly_Stats")

# --- 1) Daily_Data: append 24 rows (rows 26-49) for the new report date ---
$newDate = 46062

$newRows = @(
    @('ASAHI DEPOSITORY LLC Registered', 24895753.652, 0, 0, 0, 0, 24895753.652),
    @('ASAHI DEPOSITORY LLC Eligible', 3285306.678, 0, 628704.4, -628704.4, 0, 2656602.278),
    @('BRINK''S, INC. Registered', 17976740.579, 0, 0, 0, 0, 17976740.579),
    @('BRINK''S, INC. Eligible', 38785679.541, 0, 0, 0, 0, 38785679.541),
    @('CNT DEPOSITORY, INC. Registered', 15828675.829, 0, 0, 0, 0, 15828675.829),
    @('CNT DEPOSITORY, INC. Eligible', 13499313.638, 0, 678372.665, -678372.665, 0, 12820940.973),
    @('DELAWARE DEPOSITORY Registered', 1966294.501, 0, 0, 0, 0, 1966294.501),
    @('DELAWARE DEPOSITORY Eligible', 15756766.585, 169483.949, 19167.272, 150316.677, 0, 15907083.262),
    @('HSBC BANK, USA Registered', 3492831.93, 0, 0, 0, 0, 3492831.93),
    @('HSBC BANK, USA Eligible', 21240381.803, 0, 0, 0, 0, 21240381.803),
    @('INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Registered', 620749.47, 0, 0, 0, 0, 620749.47),
    @('INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Eligible', 3295246.644, 0, 0, 0, 0, 3295246.644),
    @('JP MORGAN CHASE BANK NA Registered', 12489589.32, 0, 0, 0, -372167.5, 12117421.82),
    @('JP MORGAN CHASE BANK NA Eligible', 155525569.983, 0, 2135547.8, -2135547.8, 372167.5, 153762189.683),
    @('LOOMIS INTERNATIONAL (US) LLC Registered', 7647170.029, 0, 0, 0, -106868.422, 7540301.607),
    @('LOOMIS INTERNATIONAL (US) LLC Eligible', 25062082.374, 0, 900673.37, -900673.37, 106868.422, 24268277.426),
    @('MALCA-AMIT ARMORED, INC. Registered', 0, 0, 0, 0, 0, 0),
    @('MALCA-AMIT ARMORED, INC. Eligible', 0, 0, 0, 0, 0, 0),
    @('MALCA-AMIT USA, LLC Registered', 1416635.864, 0, 0, 0, -191129.6, 1225506.264),
    @('MALCA-AMIT USA, LLC Eligible', 606896.577, 0, 0, 0, 191129.6, 798026.177),
    @('MANFRA, TORDELLA & BROOKES, LLC Registered', 8330589.44, 0, 0, 0, -191574.9, 8139014.54),
    @('MANFRA, TORDELLA & BROOKES, LLC Eligible', 10965243.526, 0, 0, 0, 191574.9, 11156818.426),
    @('STONEX PRECIOUS METALS LLC Registered', 7591598.24, 0, 0, 0, 0, 7591598.24),
    @('STONEX PRECIOUS METALS LLC Eligible', 186890.28, 0, 0, 0, 0, 186890.28),
)

$r = 26
foreach ($row in $newRows) {
    $wsDaily.Cells.Item($r, 1).Value = $newDate
    $wsDaily.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $wsDaily.Cells.Item($r, 2).Value = $row[0]
    $wsDaily.Cells.Item($r, 3).Value = $row[1]
    $wsDaily.Cells.Item($r, 4).Value = $row[2]
    $wsDaily.Cells.Item($r, 5).Value = $row[3]
    $wsDaily.Cells.Item($r, 6).Value = $row[4]
    $wsDaily.Cells.Item($r, 7).Value = $row[5]
    $wsDaily.Cells.Item($r, 8).Value = $row[6]
    $r = $r + 1
}

# --- 2) Today_Summary: refresh Eligible / Registered / Total_Stock for depositories whose latest totals changed ---
$summaryUpdates = @(
    @(2, 2656602.278, 24895753.652, 27552355.93),
    @(4, 12820940.973, 15828675.829, 28649616.802),
    @(5, 15907083.262, 1966294.501, 17873377.763),
    @(8, 153762189.683, 12117421.82, 165879611.503),
    @(9, 24268277.426, 7540301.607, 31808579.033),
    @(11, 798026.177, 1225506.264, 2023532.441),
    @(12, 11156818.426, 8139014.54, 19295832.966),
)
foreach ($u in $summaryUpdates) {
    $row = $u[0]
    $wsSummary.Cells.Item($row, 2).Value = $u[1]
    $wsSummary.Cells.Item($row, 3).Value = $u[2]
    $wsSummary.Cells.Item($row, 4).Value = $u[3]
}

# --- 3) Monthly_Stats: refresh the month-to-date grand total ---
$wsMonthly.Cells.Item(2, 2).Value = 284878136.493
$wsMonthly.Cells.Item(2, 3).Value = 101394888.432
$wsMonthly.Cells.Item(2, 4).Value = 386273024.925

# --- 4) Monthly_Stats: refresh cumulative RECEIVED/WITHDRAWN and latest TOTAL_TODAY per depository/region ---
$monthlyUpdates = @(
    @(7, $null, 628704.4, 2656602.278),
    @(11, $null, 2519940.473, 12820940.973),
    @(13, 169483.949, 43263.989, 15907083.262),
    @(19, $null, 2984533.7, 153762189.683),
    @(20, $null, $null, 12117421.82),
    @(21, $null, 1501353.67, 24268277.426),
    @(22, $null, $null, 7540301.607),
    @(25, $null, $null, 798026.177),
    @(26, $null, $null, 1225506.264),
    @(27, $null, $null, 11156818.426),
    @(28, $null, $null, 8139014.54),
)
foreach ($u in $monthlyUpdates) {
    $row = $u[0]
    if ($u[1] -ne $null) { $wsMonthly.Cells.Item($row, 3).Value = $u[1] }
    if ($u[2] -ne $null) { $wsMonthly.Cells.Item($row, 4).Value = $u[2] }
    if ($u[3] -ne $null) { $wsMonthly.Cells.Item($row, 5).Value = $u[3] }
}
